# "fixed dashboard and menu certificate"
# - Header in A1 was lower-case "dni"; capitalize it to "DNI".
# - The saved selection/cursor moves from the last data row (A12) up to A2
#   (first data row), as if the user clicked back into the table after
#   editing the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "dni" column header to read "DNI".
$ws.Range("A1").Value = "DNI"

# Move the active selection to A2 (top-left of the data, below the header).
$ws.Range("A2").Select()
